{"js": "// Diversion Judgment Entry: fill in the charge-result table.\n// The \"Charges\" table has key/value rows; the value cells for\n// Plea, Finding, Fine Amount, Fines Suspended and Jail Days currently\n// read \"None\" and need to be replaced with the actual adjudicated\n// values. The \"Jail Days Suspended\" row keeps its \"None\" value.\nconst replacements = [\"No Contest\", \"Guilty\", \"$ 50\", \"$ 0\", \"5\"];\n\nfor (const replacement of replacements) {\n  const results = context.document.body.search(\"None\", { matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    break;\n  }\n\n  // Replace the text of the first remaining \"None\" match in document\n  // order, in place, so the existing run formatting (font, bold, size)\n  // is preserved exactly as Word's own Find & Replace would do.\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Diversion Judgment Entry: fill in the charge-result table.\n# The \"Charges\" table has key/value rows; the value cells for\n# Plea, Finding, Fine Amount, Fines Suspended and Jail Days currently\n# read \"None\" and need to be replaced with the actual adjudicated\n# values. The \"Jail Days Suspended\" row keeps its \"None\" value.\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$tbl.Cell(4, 2).Range.Text = 'No Contest'\n$tbl.Cell(5, 2).Range.Text = 'Guilty'\n$tbl.Cell(6, 2).Range.Text = '$ 50'\n$tbl.Cell(7, 2).Range.Text = '$ 0'\n$tbl.Cell(8, 2).Range.Text = '5'\n"}
